$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    $cell.Value = ($text -replace "_old$", "_FV2404")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    $cell.Value = ($text -replace "_new$", "_FV2410")
}
